# JOHNPAUL cycle1 data fix:
#  - Every lane's CycleTime_s was re-measured/normalized to a flat 300s
#    cycle on both the Raw_Annotations (per-vehicle-type rows) and the
#    Aggregates (per-lane) sheets. Downstream throughput formulas
#    (J/K on Raw_Annotations, F/H on Aggregates) recalc automatically.
#  - Aggregates becomes the front/active sheet & selected range, whereas
#    Raw_Annotations loses its "active" marker.

$wb = $excel.ActiveWorkbook

# --- Raw_Annotations: CycleTime_s (col E, rows 2:31) -> 300 ---
$wsRaw = $wb.Worksheets.Item("Raw_Annotations")
$wsRaw.Range("E2:E31").Value = 300

# --- Aggregates: CycleTime_s (col C, rows 2:6) -> 300 ---
$wsAgg = $wb.Worksheets.Item("Aggregates")
$wsAgg.Range("C2:C6").Value = 300

# --- Selections: Raw_Annotations keeps a selection but is no longer the
#     tab shown on open; Aggregates becomes the active/selected sheet ---
[void]$wsRaw.Range("E2:E31").Select()

[void]$wsAgg.Activate()
[void]$wsAgg.Range("C2:C6").Select()
